# Artisan event sliders / commands workbook update
# - Adds a new "Artisan Command" entry: palette(<int>) / activates palette <int>
#   inserted right before the existing playbackmode(<int>) row on the
#   "Commands" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Find the row that currently holds "playbackmode(<int>)" in column B and
# insert a new row above it for the new palette(<int>) command.
$playbackRow = $ws.Cells.Find("playbackmode(<int>)").Row

$ws.Rows.Item($playbackRow).Insert()

$ws.Range("B" + $playbackRow).Value2 = "palette(<int>)"
$ws.Range("C" + $playbackRow).Value2 = "activates palette <int>"
$ws.Rows.Item($playbackRow).RowHeight = 13.8

# Update the selection/active cell to reflect the newly inserted row (mirrors
# the author's saved selection state pointing at the new palette row).
$ws.Range("B" + $playbackRow).Select()
